$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4
$ws.Range("I4").Value = "sd"
$ws.Range("J4").Value = "Statement-non-opinion"

# Row 6
$ws.Range("I6").Value = "sd"
$ws.Range("J6").Value = "Statement-non-opinion"

# Row 12
$ws.Range("I12").Value = "sd"
$ws.Range("J12").Value = "Statement-non-opinion"

# Row 17
$ws.Range("I17").Value = "sd"
$ws.Range("J17").Value = "Statement-non-opinion"

# Row 19
$ws.Range("I19").Value = "sv"
$ws.Range("J19").Value = "Statement-opinion"

# Row 25
$ws.Range("I25").Value = "sd"
$ws.Range("J25").Value = "Statement-non-opinion"

# Row 51
$ws.Range("I51").Value = "sd"
$ws.Range("J51").Value = "Statement-non-opinion"

# Row 66
$ws.Range("I66").Value = "sv"
$ws.Range("J66").Value = "Statement-opinion"
